# Update TECH_bal sheet figures for rows 4, 14, 24 (columns B:F)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TECH")

# Row 4 - Inventory
$ws.Range("B4").Value = 107000000.0
$ws.Range("C4").Value = 106000000.0
$ws.Range("D4").Value = 103000000.0
$ws.Range("E4").Value = 100000000.0
$ws.Range("F4").Value = 95000000.0

# Row 14 - Accounts Payable
$ws.Range("B14").Value = 24000000.0
$ws.Range("C14").Value = 27000000.0
$ws.Range("D14").Value = 23000000.0
$ws.Range("E14").Value = 20000000.0
$ws.Range("F14").Value = 19000000.0

# Row 24 - Long Term Tax Liability (Deferred)
$ws.Range("B24").Value = 105000000.0
$ws.Range("C24").Value = 101000000.0
$ws.Range("D24").Value = 101000000.0
$ws.Range("E24").Value = 96000000.0
$ws.Range("F24").Value = 105000000.0
